# Updated the dflowfm and rr namespaces
#
# Renames the Python module paths listed in column E ("module") of the
# "Source table" worksheet:
#   hydrolib.core.io.<x>.models      -> hydrolib.core.dflowfm.<x>.models
#   hydrolib.core.io.rr.<x>.models   -> hydrolib.core.rr.<x>.models
# (hydrolib.core.dimr.models is untouched.)

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Source table")

$map = @{
    "hydrolib.core.io.mdu.models"            = "hydrolib.core.dflowfm.mdu.models"
    "hydrolib.core.io.net.models"            = "hydrolib.core.dflowfm.net.models"
    "hydrolib.core.io.structure.models"      = "hydrolib.core.dflowfm.structure.models"
    "hydrolib.core.io.ext.models"            = "hydrolib.core.dflowfm.ext.models"
    "hydrolib.core.io.bc.models"             = "hydrolib.core.dflowfm.bc.models"
    "hydrolib.core.io.crosssection.models"   = "hydrolib.core.dflowfm.crosssection.models"
    "hydrolib.core.io.friction.models"       = "hydrolib.core.dflowfm.friction.models"
    "hydrolib.core.io.storagenode.models"    = "hydrolib.core.dflowfm.storagenode.models"
    "hydrolib.core.io.inifield.models"       = "hydrolib.core.dflowfm.inifield.models"
    "hydrolib.core.io.onedfield.models"      = "hydrolib.core.dflowfm.onedfield.models"
    "hydrolib.core.io.xyz.models"            = "hydrolib.core.dflowfm.xyz.models"
    "hydrolib.core.io.obs.models"            = "hydrolib.core.dflowfm.obs.models"
    "hydrolib.core.io.obscrosssection.models" = "hydrolib.core.dflowfm.obscrosssection.models"
    "hydrolib.core.io.rr.models"             = "hydrolib.core.rr.models"
    "hydrolib.core.io.rr.meteo.models"       = "hydrolib.core.rr.meteo.models"
    "hydrolib.core.io.rr.topology.models"    = "hydrolib.core.rr.topology.models"
}

$lastRow = 70
for ($r = 1; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 5)   # column E
    $val = $cell.Value()
    if ($val -ne $null -and $map.ContainsKey($val)) {
        $cell.Value = $map[$val]
    }
}
